# Update "Handback" report timestamps as generated by a later report run.
$wb = $excel.ActiveWorkbook

# Overview sheet: G2 "Latest HO Xliff Generate Date" for first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 07:04:15"

# zh-cn sheet: H2 "Correspond Handoff Datetime" / K2 "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 07:04:10"
$wsZhCn.Range("K2").Value = "2016-08-19 07:04:27"

# de-de sheet: H2 "Correspond Handoff Datetime" (shares text with Overview!G2)
#              K2 "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 07:04:15"
$wsDeDe.Range("K2").Value = "2016-08-19 07:04:34"
